# Finalize the import/export excel of the "user" feature:
#  - rename the template's only sheet from ADMIN to NURSE
#  - move the active selection to L5 (was I9)
#  - swap the mailto hyperlink targets that sit on A6 and A8
#  - (best-effort) restore the workbook window position/size

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet ADMIN -> NURSE
$ws.Name = "NURSE"

# 2. Update the saved selection / active cell
[void]$ws.Range("L5").Select()

# 3. Swap the hyperlink targets currently anchored on A6 and A8
#    (A6 was mailto:fdsa@gmail.com, A8 was mailto:hi@gmail.com -> swap them)
$addrRow6 = $null
$addrRow8 = $null
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 6 -and $h.Range.Column -eq 1) { $addrRow6 = $h.Address }
    elseif ($h.Range.Row -eq 8 -and $h.Range.Column -eq 1) { $addrRow8 = $h.Address }
}
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq 6 -and $h.Range.Column -eq 1) { $h.Address = $addrRow8 }
    elseif ($h.Range.Row -eq 8 -and $h.Range.Column -eq 1) { $h.Address = $addrRow6 }
}

# 4. Best-effort: restore the workbook window position/size (may be a no-op
#    in this headless runtime, harmless either way).
$win = $excel.ActiveWindow
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15720
